$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.318.39'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '3.062.53'
$ws.Range("E3").Value = '  -2.19%  '

$ws.Range("D4").Value = '''0.997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.60%  '

$ws.Range("D5").Value = '''212.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("D6").Value = '''612.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.55%  '

$ws.Range("D7").Value = '''0.365'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.12%  '

$ws.Range("D8").Value = '''0.890'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +24.61%  '

$ws.Range("D9").Value = '''0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.24%  '

$ws.Range("D10").Value = '3.060.13'
$ws.Range("E10").Value = '  -2.15%  '

$ws.Range("D11").Value = '''0.673'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +22.58%  '

$ws.Range("D12").Value = '''0.188'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.49%  '

$ws.Range("D13").Value = '''0.0000238'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.63%  '

$ws.Range("D14").Value = '''5.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.46%  '

$ws.Range("D15").Value = '89.460.88'
$ws.Range("E15").Value = '  -0.38%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '''32.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.30%  '

$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '3.627.14'
$ws.Range("E17").Value = '  -2.88%  '

$ws.Range("D18").Value = '3.126.84'
$ws.Range("E18").Value = '  -2.66%  '

$ws.Range("D19").Value = '''3.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.35%  '

$ws.Range("D20").Value = '''0.0000211'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.12%  '

$ws.Range("D21").Value = '''13.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.78%  '

$ws.Range("D22").Value = '''427.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.28%  '

$ws.Range("D23").Value = '''4.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.77%  '

$ws.Range("D24").Value = '''8.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.64%  '

$ws.Range("D25").Value = '''5.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.90%  '

$ws.Range("D26").Value = '''84.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.84%  '

$ws.Range("D27").Value = '''11.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.13%  '

$ws.Range("D28").Value = '3.237.94'
$ws.Range("E28").Value = '  -3.40%  '

$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E30").Value = '  +9.20%  '

$ws.Range("D31").Value = '''0.162'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.55%  '

$ws.Range("D32").Value = '''8.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.64%  '

$ws.Range("D33").Value = '''504.59'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("D34").Value = '''3.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.78%  '

$ws.Range("D35").Value = '''6.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.43%  '

$ws.Range("D36").Value = '''22.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.84%  '

$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '''1.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.83%  '

$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").Value = '''1.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.69%  '

$ws.Range("E39").Value = '  +5.49%  '

$ws.Range("E40").Value = '  -0.22%  '

$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.39%  '

$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("E43").Value = '  +13.46%  '

$ws.Range("D44").Value = '''0.367'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.19%  '

$ws.Range("D45").Value = '''1.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.39%  '

$ws.Range("D46").Value = '''147.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("D47").Value = '''0.0695'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +15.57%  '

$ws.Range("D48").Value = '''43.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.89%  '

$ws.Range("D49").Value = '''4.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.47%  '

$ws.Range("D50").Value = '''1.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.61%  '

$ws.Range("D51").Value = '''155.90'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.14%  '

